$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the report title (row 1) for the new RTO / state.
$ws.Range("A1").Value = "Maker Month Wise Data  of SAHEBGANJ - JH18 , Jharkhand (2022)"

# 2) Add the new data row (row 5) describing OLA ELECTRIC TECHNOLOGIES PVT LTD.
#    Month columns: C=JAN D=FEB E=MAR F=APR G=MAY H=JUN I=JUL J=AUG K=SEP L=OCT M=NOV N=DEC O=TOTAL
$rowValues = @{
    "A5" = "1"
    "B5" = "OLA ELECTRIC TECHNOLOGIES PVT LTD"
    "C5" = "0"
    "D5" = "0"
    "E5" = "0"
    "F5" = "0"
    "G5" = "0"
    "H5" = "2"
    "I5" = "0"
    "J5" = "0"
    "K5" = "0"
    "L5" = "1"
    "M5" = "0"
    "N5" = "0"
    "O5" = "3"
}

foreach ($addr in $rowValues.Keys) {
    $val = $rowValues[$addr]
    $cell = $ws.Range($addr)
    # Force numeric-looking text ("0","1","2","3") to be stored as text (shared string)
    # rather than being auto-converted to a number, matching the source data which
    # keeps these as text cells. Non-numeric text (maker name) is unaffected by this.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# 3) Keep a trailing blank row (was row 9, now row 10) and extend the sheet dimension
#    down to row 10, matching the original template's trailing blank row.
$ws.Range("A10").Font.Size = 11

# 4) Resize column A/O (the thin spacer columns) and column B (maker name column)
#    so the new, much longer maker name fits / the layout matches the updated sheet.
$ws.Columns.Item(1).ColumnWidth = 1.3
$ws.Columns.Item(2).ColumnWidth = 36.6
$ws.Columns.Item(15).ColumnWidth = 1.3

Write-Output "Applied Jharkhand SAHEBGANJ - JH18 data edits"
